$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.752.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.05%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.277.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.22%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'251.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.643"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.59%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'74.81"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +6.42%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.644"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'39.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.45%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0979"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.72%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.22%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.93%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.618.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.48%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.89%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.866"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.48%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.267.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.17%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.640.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.18%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.13%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'72.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'236.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.62%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +6.30%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.82%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D26").Value = "'11.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.22%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.82%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.32%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'167.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.03%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'21.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.39%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0874"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +8.84%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.55%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.126"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.53%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'31.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.71%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.128"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.95%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.70%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.61%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -5.17%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'13.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +9.34%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.10%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +3.51%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'61.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.34%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.35%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'105.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +11.44%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'4.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.60%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.80%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.25%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.15%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.42%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'4.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.71%  "
$ws.Range("E51").Style = "Normal"
